$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring formatting for the new row in line with the existing data rows by
# copying row 7's formats down into row 8 (values are overwritten next).
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)  # xlPasteFormats

# New timesheet entry: 11/11/2020, 3:49 PM - 5:56 PM, 10 min interrupt,
# "writing requirements/design documents".
$ws.Range("A8").Value = 44146
$ws.Range("B8").Value = 0.65902777777777777
$ws.Range("C8").Value = 0.74722222222222223
$ws.Range("D8").Value = 10
$ws.Range("F8").Value = "writing requirements/design documents"

# Extend the Delta formula down through the new row; setting the formula on
# the whole range at once lets Excel share the formula across E2:E8.
$ws.Range("E2:E8").Formula = "=C2-B2-TIME(0,D2,0)"

# Matches the author's final selection before saving.
[void]$ws.Range("C9").Select()
